# Applies the "gpu-cores" update:
#  - A4: "Execution Unit" -> "Execution Unit/per Core"
#  - A5: "ALU" -> "ALU/per EU"
#  - C3: numeric 3 -> text "unknow"
#  - C5: "unknow" -> "unknow（Maybe Total 1024ALU）"
#  - A8: new summary text (Chinese)
#  - Selection moves to A5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Execution Unit/per Core"
$ws.Range("A5").Value = "ALU/per EU"

$ws.Range("C3").Value = "unknow"
$ws.Range("C5").Value = "unknow（Maybe Total 1024ALU）"

$ws.Range("A8").Value = "综合来看，Apple GPU的频率较高，较多的EU能够提升硬件大粒度调度的灵活性，算力最强。"

[void]$ws.Range("A5").Select()
